$d = $word.ActiveDocument

# --- Change 1: remove the stray "_GoBack" bookmark that currently sits
#     on the very first paragraph ("Микроконтроллеры:") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2: delete the whole "Гироскоп:" block (the heading paragraph
#     plus its three bullet items) and move the "_GoBack" bookmark to the
#     start of the paragraph that follows it ("Ик датчик:") ---

# Locate the start of the "Гироскоп:" heading paragraph (the heading text
# is the very first thing in that paragraph, so the match start IS the
# paragraph start).
$headingRange = $d.Content
$headingRange.Find.ClearFormatting()
$headingFound = $headingRange.Find.Execute("Гироскоп:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $headingFound) {
    throw "Could not find the 'Гироскоп:' heading paragraph"
}
$headingParaStart = $headingRange.Start

# Locate the start of the following "Ик датчик:" paragraph the same way;
# this is also where the deletion must stop.
$nextRange = $d.Content
$nextRange.Find.ClearFormatting()
$nextFound = $nextRange.Find.Execute("Ик датчик:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $nextFound) {
    throw "Could not find the 'Ик датчик:' paragraph"
}
$nextParaStart = $nextRange.Start

$gyroRange = $d.Range($headingParaStart, $nextParaStart)
$gyroRange.Delete()

# Re-insert the "_GoBack" bookmark right at the start of the (now
# immediately following) "Ик датчик:" paragraph.
$irRange = $d.Content
$irRange.Find.ClearFormatting()
$irRange.Find.Execute("Ик датчик:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$irParaStart = $irRange.Start
$bmRange = $d.Range($irParaStart, $irParaStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Change 3: drop the stray <w:lastRenderedPageBreak/> recorded on the
#     very last paragraph, without touching its run's text/formatting ---
$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$lastStart = $lastParagraph.Range.Start
$insertionPoint = $d.Range($lastStart, $lastStart)
$insertionPoint.InsertAfter("X")
$tempCharRange = $d.Range($lastStart, $lastStart + 1)
$tempCharRange.Delete()

Write-Output "Done"
